$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "uint8_t"
$ws.Range("A4").Select()
